# Commit: "Fruta / hortaliza, semanal"
# Weekly data refresh: a new pricing group (3 rows, one per "Calidad":
# Especial / Primera / Segunda) for Femacal de La Calera - Piña is added
# for the date 2021-11-11 (Excel serial 44511). It is inserted right
# before the existing row 308, which pushes every following row down by
# three positions (old row N now lives at row N+3). The sheet's
# dimension therefore grows from A1:T386 to A1:T389.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 308, shifting the rest
# of the table (rows 308..386) down to 311..389.
$ws.Range("A308:A310").EntireRow.Insert()

# Shared constant values for every data row in this sheet.
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$codreg = 5
$tipo = "Fruta"
$productoId = 100108
$producto = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria = "Piña"
$variedad = "Caramelo"
$origen = "Ecuador"
$fecha = 44511

# The three new rows (Especial / Primera / Segunda) for the new date.
$newRows = @(
    @{ Row = 308; Calidad = "Especial"; Volumen = 108; Min = 19000; Max = 19000; Prom = 19000; Unidad = "`$/caja 10 unidades"; PrecioKg = 1900; KgUnidad = 10 },
    @{ Row = 309; Calidad = "Primera";  Volumen = 162; Min = 19000; Max = 19000; Prom = 19000; Unidad = "`$/caja 12 unidades"; PrecioKg = 1583; KgUnidad = 12 },
    @{ Row = 310; Calidad = "Segunda";  Volumen = 162; Min = 19000; Max = 19000; Prom = 19000; Unidad = "`$/caja 14 unidades"; PrecioKg = 1357; KgUnidad = 14 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
